# console/app_cover cover-slide refresh: lab name, version/date, and the
# tagline text+layout on the first ("cover") slide.
#
# EMU<->point note: Shape.Left/Top/Width/Height are exposed in points
# (1 pt = 12700 EMU) while the underlying OOXML stores EMU integers. The
# runtime truncates when converting points back to EMU, so a bare
# `emu / 12700.0` can land one EMU short through float error. Nudging by
# half an EMU before dividing removes that bias and round-trips exactly.
function EmuToPt([double]$emu) {
    return ($emu + 0.5) / 12700.0
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- Shape id=13 "TextBox 12": lab name -------------------------------
$shLab = $s.Shapes.Item(4)
$shLab.TextFrame.TextRange.Text = "Surgical Robotics and Instrumentation Lab"

# --- Shape id=15 "TextBox 14": version + date -------------------------
$shVer = $s.Shapes.Item(6)
$trVer = $shVer.TextFrame.TextRange
$origVerHeight = $shVer.Height

$txt = $trVer.Text
$i = $txt.IndexOf("V0.2.1")
$trVer.Characters($i + 1, 6).Text = "V0.3.0"

$txt = $trVer.Text
$i = $txt.IndexOf("September 20, 2021")
$trVer.Characters($i + 1, 19).Text = "August 2023"

# Editing the field runs nudges the shape's auto-fit height; this box's
# size is untouched in the target, so restore it precisely.
$shVer.Height = $origVerHeight

# --- Shape id=16 "TextBox 15": tagline text, centering, resize/move ---
$shTag = $s.Shapes.Item(7)
$trTag = $shTag.TextFrame.TextRange
$trTag.ParagraphFormat.Alignment = 2   # ppAlignCenter
$trTag.Text = "Configure and control the MRI compatible robotic system"
$shTag.Left = EmuToPt 705181
$shTag.Height = EmuToPt 446020
